$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.983.98'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.64%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.522.26'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.84%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.17%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '605.10'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.09%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.98'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.00%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.519.89'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.80%  '

$ws.Range("E8").Value = '  +0.17%  '

$ws.Range("E9").Value = '  -1.16%  '

$ws.Range("E10").Value = '  +0.16%  '

$ws.Range("E11").Value = '  +3.36%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.424'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.44%  '

$ws.Range("E13").Value = '  +0.16%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.116.82'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.83%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '31.77'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.57%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.510.45'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.28%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.983.04'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.53%  '

$ws.Range("E18").Value = '  -0.02%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.82'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +10.31%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.40'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.74%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.39'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.88%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '436.55'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.80%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.612'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.87%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '79.62'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.10%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.659.30'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.69%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000122'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.66%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.87'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.60%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.43'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.82%  '

$ws.Range("E30").Value = '  +0.44%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.59'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.53%  '

$ws.Range("E32").Value = '  +0.53%  '

$ws.Range("E33").Value = '  -0.63%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '25.45'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.56%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.511.31'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.78%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.81'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.61%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.90'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.78%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.07'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.19%  '

$ws.Range("E39").Value = '  +0.00%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.998'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.16%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0894'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.71%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '169.56'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.17%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.45'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.38%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.09'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -9.74%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.897'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.21%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '45.91'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.27%  '

$ws.Range("B47").Value = 'ONDO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.32'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.36%  '

$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '28.35'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.59%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.50'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.37%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.45'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.10%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.997'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.79%  '
